# Fix researcher name typo: "Xianfeng Zhang" -> "Xianfeng Zeng"
# The name appears in the "Researcher Name" column (C) of the Samples sheet,
# for every sample row (C2:C17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")
$ws.Range("C2:C17").Value = "Xianfeng Zeng"
